# Weekly fruit/vegetable price update: insert a new weekly record as the
# new row 75 (Fecha 45167 = 2023-08-29) for "Feria Lagunitas de Puerto
# Montt" / Membrillo / Champion, shifting the existing rows 75-179 down to
# 76-180.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 75, pushing everything below
# down by one (old row 75 -> 76, ..., old row 179 -> 180).
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new weekly observation.
$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C75").Value = "Los Lagos"
$ws.Range("D75").Value = 45167
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100104
$ws.Range("H75").Value = "Frutos de pepita"
$ws.Range("I75").Value = 100104003
$ws.Range("J75").Value = "Membrillo"
$ws.Range("K75").Value = "Champion"
$ws.Range("L75").Value = "Primera"
$ws.Range("M75").Value = 200
$ws.Range("N75").Value = 15000
$ws.Range("O75").Value = 15000
$ws.Range("P75").Value = 15000
$ws.Range("Q75").Value = "$/caja 18 kilos empedrada"
$ws.Range("R75").Value = "Región de O'Higgins"
$ws.Range("S75").Value = 833
$ws.Range("T75").Value = 18
